# Update heterogeneity meta-analysis precipitation table
# to also include PdeltaAIC as a cov for the CG path.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell while preserving the cell's
# original (unstyled / "General") formatting. The target cells hold
# numeric-looking strings (e.g. " 195.9") that must remain stored as
# text (as in the source workbook), not be auto-converted to numbers.
# Trick: temporarily borrow a "text" number format (@) from a cell
# that already uses it (F2), assign the plain value (now interpreted
# as text because the cell format is already "@"), then restore the
# original default styling by pasting formats from an unstyled cell
# (C4, Trait_Category column) so no stray style entries are left
# behind.
function Set-PlainTextValue($addr, $val) {
    $ws.Range("F2").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = $val
    $ws.Range("C4").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# Column F cells already use the "text" (@) number format in the
# original workbook (same style as F2), so just (re)assign the value
# after ensuring that style is in place - no restyling needed after.
function Set-PvalTextValue($addr, $val) {
    $ws.Range("F2").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = $val
}

# Row 4: Morphological / Precipitation / CG
Set-PlainTextValue "A4" " 195.9"
Set-PlainTextValue "B4" "0.429"

# Row 6: Morphological / Precipitation / CZG
Set-PlainTextValue "A6" " 112.3"
Set-PvalTextValue "F6" "0.473"

# Row 7: Morphological / Precipitation / TotalCG
Set-PlainTextValue "A7" " 153.6"
Set-PlainTextValue "B7" "0.351"
Set-PvalTextValue "F7" "0.006"

# Row 10: Phenological / Precipitation / CG
Set-PlainTextValue "A10" " 136.1"
Set-PlainTextValue "B10" "0.415"

# Row 12: Phenological / Precipitation / CZG
Set-PlainTextValue "A12" " 101.9"
Set-PlainTextValue "B12" "0.059"
Set-PvalTextValue "F12" "0.224"

# Row 13: Phenological / Precipitation / TotalCG
Set-PlainTextValue "A13" "  58.6"
